$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 470.2
$ws.Range("J43").Value = 487.5
$ws.Range("L43").Value = 487.5
$ws.Range("N43").Value = -625.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 473.5
$ws.Range("I4").Value = 445
$ws.Range("J4").Value = 502
$ws.Range("K4").Value = 445
$ws.Range("L4").Value = 502
$ws.Range("M4").Value = -329
$ws.Range("N4").Value = -734

$ws.Range("H10").Value = 502
$ws.Range("I10").Value = 502
$ws.Range("K10").Value = 502
$ws.Range("M10").Value = -332

$ws.Range("H22").Value = 15254
$ws.Range("I22").Value = 15254
$ws.Range("K22").Value = 15254
$ws.Range("M22").Value = -14955

$ws.Range("H32").Value = 12503021
$ws.Range("I32").Value = 3452.8572
$ws.Range("J32").Value = 100000000
$ws.Range("K32").Value = 3452.8572
$ws.Range("L32").Value = 100000000
$ws.Range("M32").Value = -3165.8572
$ws.Range("N32").Value = -100000574

$ws.Range("H35").Value = 5317.5713
$ws.Range("I35").Value = 2394.6
$ws.Range("K35").Value = 2394.6
$ws.Range("M35").Value = -1988.6

$ws.Range("H110").Value = 3149.6667
$ws.Range("I110").Value = 2999.5
$ws.Range("K110").Value = 2999.5
$ws.Range("M110").Value = -954.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 767.375
$ws.Range("I29").Value = 767.375
$ws.Range("K29").Value = 767.375
$ws.Range("M29").Value = -478.375

$ws.Range("H36").Value = 989.5714
$ws.Range("I36").Value = 989.5714
$ws.Range("K36").Value = 989.5714
$ws.Range("M36").Value = -455.5714

$ws.Range("H86").Value = 4952.778
$ws.Range("I86").Value = 2679.1667
$ws.Range("J86").Value = 9500
$ws.Range("K86").Value = 2679.1667
$ws.Range("L86").Value = 9500
$ws.Range("M86").Value = -1556.1667
$ws.Range("N86").Value = -11746

$ws.Range("H89").Value = 4952.778
$ws.Range("I89").Value = 2679.1667
$ws.Range("J89").Value = 9500
$ws.Range("K89").Value = 13395.8335
$ws.Range("L89").Value = 47500
$ws.Range("M89").Value = -7779.833500000001
$ws.Range("N89").Value = -58732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.36
$ws.Range("I7").Value = 47.416668
$ws.Range("J7").Value = 89.61539
$ws.Range("K7").Value = 47.416668
$ws.Range("L7").Value = 89.61539
$ws.Range("M7").Value = 65.583332
$ws.Range("N7").Value = -315.61539

$ws.Range("H22").Value = 1291.7142
$ws.Range("I22").Value = 324
$ws.Range("J22").Value = 2582
$ws.Range("K22").Value = 324
$ws.Range("L22").Value = 2582
$ws.Range("M22").Value = 26
$ws.Range("N22").Value = -3282

$ws.Range("H31").Value = 5482.0645
$ws.Range("I31").Value = 3704.7
$ws.Range("K31").Value = 3704.7
$ws.Range("M31").Value = -3409.7

$ws.Range("H34").Value = 5482.0645
$ws.Range("I34").Value = 3704.7
$ws.Range("K34").Value = 3704.7
$ws.Range("M34").Value = -3502.7

$ws.Range("H42").Value = 13200
$ws.Range("I42").Value = 13200
$ws.Range("K42").Value = 13200
$ws.Range("M42").Value = -12607

$ws.Range("H55").Value = 17255
$ws.Range("J55").Value = 18882.5
$ws.Range("L55").Value = 18882.5
$ws.Range("N55").Value = -19512.5

$ws.Range("H62").Value = 2333.3333
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2333.3333
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21240

$ws.Range("H70").Value = 35892.855
$ws.Range("J70").Value = 35892.855
$ws.Range("L70").Value = 35892.855
$ws.Range("N70").Value = -36522.855

$ws.Range("H73").Value = 35892.855
$ws.Range("J73").Value = 35892.855
$ws.Range("L73").Value = 35892.855
$ws.Range("N73").Value = -38076.855

$ws.Range("H97").Value = 58000
$ws.Range("J97").Value = 58000
$ws.Range("L97").Value = 58000
$ws.Range("N97").Value = -59982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31.4375
$ws.Range("I2").Value = 21.666666
$ws.Range("J2").Value = 44
$ws.Range("K2").Value = 129.999996
$ws.Range("L2").Value = 264
$ws.Range("M2").Value = -16.99999600000001
$ws.Range("N2").Value = -490

$ws.Range("H33").Value = 252.44444
$ws.Range("J33").Value = 85.666664
$ws.Range("L33").Value = 513.999984
$ws.Range("N33").Value = -1079.999984

$ws.Range("H49").Value = 2992.5
$ws.Range("I49").Value = 2985
$ws.Range("K49").Value = 8955
$ws.Range("M49").Value = -8799

$ws.Range("H80").Value = 4192.8
$ws.Range("I80").Value = 3959.95
$ws.Range("K80").Value = 11879.85
$ws.Range("M80").Value = -10943.85

$ws.Range("H83").Value = 4192.8
$ws.Range("I83").Value = 3959.95
$ws.Range("K83").Value = 35639.55
$ws.Range("M83").Value = -30959.55

$ws.Range("H107").Value = 446.47058
$ws.Range("I107").Value = 229.29411
$ws.Range("J107").Value = 663.64703
$ws.Range("K107").Value = 687.8823299999999
$ws.Range("L107").Value = 1990.94109
$ws.Range("M107").Value = 1232.11767
$ws.Range("N107").Value = -5830.94109

$ws.Range("H114").Value = 1314.3077
$ws.Range("J114").Value = 1180.125
$ws.Range("L114").Value = 3540.375
$ws.Range("N114").Value = -10048.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 311.42105
$ws.Range("I2").Value = 152.875
$ws.Range("J2").Value = 426.72726
$ws.Range("K2").Value = 152.875
$ws.Range("L2").Value = 426.72726
$ws.Range("M2").Value = -39.875
$ws.Range("N2").Value = -652.72726

$ws.Range("H3").Value = 24302046
$ws.Range("I3").Value = 17139108
$ws.Range("J3").Value = 44000130
$ws.Range("K3").Value = 17139108
$ws.Range("L3").Value = 44000130
$ws.Range("M3").Value = -17138992
$ws.Range("N3").Value = -44000362

$ws.Range("H102").Value = 2035.5
$ws.Range("I102").Value = 1434
$ws.Range("K102").Value = 1434
$ws.Range("M102").Value = 188

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2061.353
$ws.Range("I16").Value = 2103.4
$ws.Range("K16").Value = 2103.4
$ws.Range("M16").Value = -1933.4

$ws.Range("H17").Value = 5166.6665
$ws.Range("I17").Value = 5500
$ws.Range("J17").Value = 4500
$ws.Range("K17").Value = 5500
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = -5330
$ws.Range("N17").Value = -4840

$ws.Range("H40").Value = 5628.5454
$ws.Range("I40").Value = 5322.625
$ws.Range("J40").Value = 6444.3335
$ws.Range("K40").Value = 5322.625
$ws.Range("L40").Value = 6444.3335
$ws.Range("M40").Value = -5186.625
$ws.Range("N40").Value = -6716.3335

$ws.Range("H46").Value = 6296.7334
$ws.Range("I46").Value = 2266.6667
$ws.Range("J46").Value = 7304.25
$ws.Range("K46").Value = 2266.6667
$ws.Range("L46").Value = 7304.25
$ws.Range("M46").Value = -2078.6667
$ws.Range("N46").Value = -7680.25

$ws.Range("H55").Value = 1398.75
$ws.Range("I55").Value = 1255
$ws.Range("J55").Value = 1686.25
$ws.Range("K55").Value = 1255
$ws.Range("L55").Value = 1686.25
$ws.Range("M55").Value = -1082
$ws.Range("N55").Value = -2032.25

$ws.Range("H61").Value = 6424.25
$ws.Range("I61").Value = 4973.5
$ws.Range("K61").Value = 4973.5
$ws.Range("M61").Value = -4771.5

$ws.Range("H100").Value = 7374.875
$ws.Range("I100").Value = 4833
$ws.Range("J100").Value = 8900
$ws.Range("K100").Value = 4833
$ws.Range("L100").Value = 8900
$ws.Range("M100").Value = -4292
$ws.Range("N100").Value = -9982

$ws.Range("H113").Value = 6424.25
$ws.Range("I113").Value = 4973.5
$ws.Range("K113").Value = 4973.5
$ws.Range("M113").Value = -2803.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3100
$ws.Range("I17").Value = 3500
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 3500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = -3328
$ws.Range("N17").Value = -2844

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
